$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to remain a text string (so that
# Excel does not auto-convert number-looking strings like "505.09" into a
# numeric value), and then restore the cell's style so no stray number-format
# style gets left behind on the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Simple price / volume(1h) updates (rows where only D and/or E changed) ---
$ws.Range("D2").Value = "60.009.79"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.582.06"
$ws.Range("E3").Value = "  -0.29%  "
Set-TextValue "D5" "505.09"
$ws.Range("E5").Value = "  +0.00%  "
Set-TextValue "D6" "152.85"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -7.74%  "
$ws.Range("D9").Value = "2.591.36"
$ws.Range("E9").Value = "  +0.27%  "
Set-TextValue "D10" "6.65"
$ws.Range("E10").Value = "  +7.70%  "
$ws.Range("E11").Value = "  -0.20%  "
Set-TextValue "D12" "0.345"
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "3.043.60"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "60.181.83"
$ws.Range("E15").Value = "  +1.39%  "
Set-TextValue "D16" "21.49"
$ws.Range("E16").Value = "  -1.43%  "
Set-TextValue "D17" "0.0000139"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "2.595.83"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +1.97%  "
Set-TextValue "D20" "345.66"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  -1.04%  "
Set-TextValue "D24" "60.24"
$ws.Range("E24").Value = "  +0.37%  "
Set-TextValue "D25" "0.420"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("D27").Value = "2.701.31"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +2.55%  "
Set-TextValue "D30" "7.38"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -0.06%  "
Set-TextValue "D32" "19.29"
$ws.Range("E32").Value = "  -0.58%  "
Set-TextValue "D33" "153.40"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("E34").Value = "  -0.82%  "
Set-TextValue "D35" "5.74"
$ws.Range("E35").Value = "  +4.86%  "
Set-TextValue "D36" "3.98"
$ws.Range("E36").Value = "  +0.90%  "
Set-TextValue "D37" "1.18"
$ws.Range("E37").Value = "  -0.67%  "
Set-TextValue "D42" "35.83"
$ws.Range("E42").Value = "  +2.28%  "
Set-TextValue "D43" "295.83"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  -1.66%  "
Set-TextValue "D48" "19.72"
$ws.Range("E48").Value = "  +1.09%  "
Set-TextValue "D49" "4.87"
$ws.Range("E49").Value = "  +2.18%  "
Set-TextValue "D50" "0.0233"
$ws.Range("E50").Value = "  -1.90%  "
Set-TextValue "D51" "10.28"
$ws.Range("E51").Value = "  -0.19%  "

# --- Rows that were reordered (coin identity swapped between two adjacent rows) ---
# Row 38 (was Fetch.AI) -> now SuiNetwork
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D38" "0.842"
$ws.Range("E38").Value = "  +17.42%  "

# Row 39 (was SuiNetwork) -> now Fetch.AI
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D39" "0.844"
$ws.Range("E39").Value = "  -1.25%  "

# Row 40 (was Filecoin) -> now Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "1.46"
$ws.Range("E40").Value = "  +1.86%  "

# Row 41 (was Stacks) -> now Filecoin
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D41" "3.75"
$ws.Range("E41").Value = "  -0.45%  "

# Row 46 (was Hedera) -> now FirstDigitalUSD
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D46" "0.999"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47 (was FirstDigitalUSD) -> now Hedera
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0557"
$ws.Range("E47").Value = "  -1.09%  "
